$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string needed for the "Neutrophils" cluster (inserted between
# "Inflammatory-Mac" and "Resolving-Mac" in the sending/target cluster set).

# Full target table: cross join of the three clusters (Inflammatory-Mac,
# Neutrophils, Resolving-Mac) as Sending x Target cluster, with ligand
# Ccl24 / receptor Ccr3 and the refreshed TPM-derived metrics.

$rows = @(
  @{ Row=2;  A="Inflammatory-Mac"; D="Inflammatory-Mac"; E=3; F=1; G=2.208740666666666; H=6.626221999999999; I=0.5931493987290414; J=0.5931493987290414; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846;   O=0.1293260700537641; P=0.1293260700537641; Q=0.3591750997568889;  R=3.232575897812;    S=0.07670968069238009; T=0.07670968069238009 },
  @{ Row=3;  A="Inflammatory-Mac"; D="Neutrophils";      E=3; F=1; G=2.208740666666666; H=6.626221999999999; I=0.5931493987290414; J=0.5931493987290414; K=3; L=1;                M=0.8767803333333334; N=2.630341;    O=0.6972931302732585; P=0.6972931302732585; Q=1.936580377966889;  R=17.429223401702;   S=0.4135990009594744;  T=0.4135990009594744  },
  @{ Row=4;  A="Inflammatory-Mac"; D="Resolving-Mac";    E=3; F=1; G=2.208740666666666; H=6.626221999999999; I=0.5931493987290414; J=0.5931493987290414; K=3; L=1;                M=0.21801;            N=0.65403;    O=0.1733807996729775; P=0.1733807996729775; Q=0.48152755274;      R=4.33374797466;     S=0.1028407170771869;  T=0.1028407170771869  },
  @{ Row=5;  A="Neutrophils";      D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.1917523333333333; H=0.575257; I=0.05149440264221032; J=0.05149440264221032; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846;   O=0.1293260700537641; P=0.1293260700537641; Q=0.03118186960244444; R=0.280636826422;    S=0.006659568723483229; T=0.006659568723483229 },
  @{ Row=6;  A="Neutrophils";      D="Neutrophils";      E=1; F=0.3333333333333333; G=0.1917523333333333; H=0.575257; I=0.05149440264221032; J=0.05149440264221032; K=3; L=1;                M=0.8767803333333334; N=2.630341;    O=0.6972931302732585; P=0.6972931302732585; Q=0.1681246747374444;  R=1.513122072637;    S=0.03590669320993838;  T=0.03590669320993838  },
  @{ Row=7;  A="Neutrophils";      D="Resolving-Mac";    E=1; F=0.3333333333333333; G=0.1917523333333333; H=0.575257; I=0.05149440264221032; J=0.05149440264221032; K=3; L=1;                M=0.21801;            N=0.65403;    O=0.1733807996729775; P=0.1733807996729775; Q=0.04180392619;       R=0.37623533571;     S=0.008928140708788709; T=0.008928140708788709 },
  @{ Row=8;  A="Resolving-Mac";    D="Inflammatory-Mac"; E=3; F=1; G=1.323258;           H=3.969774; I=0.3553561986287483;  J=0.3553561986287483;  K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846;   O=0.1293260700537641; P=0.1293260700537641; Q=0.215182040756;      R=1.936638366804;    S=0.04595682063790082;  T=0.04595682063790082  },
  @{ Row=9;  A="Resolving-Mac";    D="Neutrophils";      E=3; F=1; G=1.323258;           H=3.969774; I=0.3553561986287483;  J=0.3553561986287483;  K=3; L=1;                M=0.8767803333333334; N=2.630341;    O=0.6972931302732585; P=0.6972931302732585; Q=1.160206590326;      R=10.441859312934;   S=0.2477874361038457;  T=0.2477874361038457  },
  @{ Row=10; A="Resolving-Mac";    D="Resolving-Mac";    E=3; F=1; G=1.323258;           H=3.969774; I=0.3553561986287483;  J=0.3553561986287483;  K=3; L=1;                M=0.21801;            N=0.65403;    O=0.1733807996729775; P=0.1733807996729775; Q=0.28848347658;       R=2.59635128922;     S=0.06161194188700179;  T=0.06161194188700179  }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value  = $r.A
  $ws.Cells.Item($row, 2).Value  = "Ccl24"
  $ws.Cells.Item($row, 3).Value  = "Ccr3"
  $ws.Cells.Item($row, 4).Value  = $r.D
  $ws.Cells.Item($row, 5).Value  = $r.E
  $ws.Cells.Item($row, 6).Value  = $r.F
  $ws.Cells.Item($row, 7).Value  = $r.G
  $ws.Cells.Item($row, 8).Value  = $r.H
  $ws.Cells.Item($row, 9).Value  = $r.I
  $ws.Cells.Item($row, 10).Value = $r.J
  $ws.Cells.Item($row, 11).Value = $r.K
  $ws.Cells.Item($row, 12).Value = $r.L
  $ws.Cells.Item($row, 13).Value = $r.M
  $ws.Cells.Item($row, 14).Value = $r.N
  $ws.Cells.Item($row, 15).Value = $r.O
  $ws.Cells.Item($row, 16).Value = $r.P
  $ws.Cells.Item($row, 17).Value = $r.Q
  $ws.Cells.Item($row, 18).Value = $r.R
  $ws.Cells.Item($row, 19).Value = $r.S
  $ws.Cells.Item($row, 20).Value = $r.T
}
